# Update the "Estado de Cuenta" data:
#  - Salario Basico (column G) for the three detail rows moves from
#    4,367,944 to 4,492,340.
#  - The "Periodo Mora" labels (column E) for rows 16 and 18 swap, so the
#    first detail row now reads "2011" and the last reads "2101" (the
#    middle row stays "2012").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2011"
$ws.Range("E18").Value = "2101"

$ws.Range("G16").Value = 4492340
$ws.Range("G17").Value = 4492340
$ws.Range("G18").Value = 4492340
